$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper values (new localization "handoff" rows being reported for this run)
# ---------------------------------------------------------------------------
$readyForHandoff   = "Ready for handoff"
$overviewDate      = "2016-03-24 15:23:12"
$zhDate            = "2016-03-24 15:23:05"
$epoch             = "0001-01-01 00:00:00"

$file1             = "51d2b0be-b8f8-4d6b-b27a-a5b90c660bdc.png"
$file2             = "921bfe5e-515b-4217-b7f1-b24b8e912cb7.png"
$file3             = "e2efc184-d807-410b-ac76-495012fc08d8.md"

$file1ext          = ".png"
$file3ext          = ".md"

$target1           = "b8783db50389cae3be61bbeda01d94acc38894ff.png"
$target2           = "574a8378940309a39647bff10e071574a08f1c05.png"
$target3zh         = "e2efc184-d807-410b-ac76-495012fc08d8.08518101bfe08f50e3652bae0fd09588ed537b98.zh-cn.xlf"
$target3de         = "e2efc184-d807-410b-ac76-495012fc08d8.08518101bfe08f50e3652bae0fd09588ed537b98.de-de.xlf"

$isDependency      = "IsDependency"
$include           = "Include"
$dependencyFrom    = "e2e\e2efc184-d807-410b-ac76-495012fc08d8.md"

# ===========================================================================
# Sheet 1: Overview
# ===========================================================================
$ws1 = $wb.Worksheets.Item("Overview")

# Duplicate row 2 (with its formatting/hyperlink-style) into rows 3 and 4
$ws1.Rows(2).Copy()
$ws1.Rows(3).Insert(-4121)
$ws1.Rows(2).Copy()
$ws1.Rows(4).Insert(-4121)

# Row 2: existing handoff file renamed/updated
$ws1.Range("A2").Value2 = $file1
$ws1.Range("B2").Value2 = $readyForHandoff
$ws1.Range("C2").Value2 = $readyForHandoff
$ws1.Range("D2").Value2 = $overviewDate

# Row 3: new handoff file
$ws1.Range("A3").Value2 = $file2
$ws1.Range("B3").Value2 = $readyForHandoff
$ws1.Range("C3").Value2 = $readyForHandoff
$ws1.Range("D3").Value2 = $overviewDate

# Row 4: new handoff file
$ws1.Range("A4").Value2 = $file3
$ws1.Range("B4").Value2 = $readyForHandoff
$ws1.Range("C4").Value2 = $readyForHandoff
$ws1.Range("D4").Value2 = $overviewDate

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$file1", "", "", $file1)
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$file2", "", "", $file2)
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$file3", "", "", $file3)

# ===========================================================================
# Sheet 2: zh-cn
# ===========================================================================
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows(2).Copy()
$ws2.Rows(3).Insert(-4121)
$ws2.Rows(2).Copy()
$ws2.Rows(4).Insert(-4121)

# Row 2
$ws2.Range("A2").Value2 = $file1
$ws2.Range("B2").Value2 = $file1ext
$ws2.Range("C2").Value2 = $readyForHandoff
$ws2.Range("D2").Value2 = $target1
$ws2.Range("E2").Value2 = $zhDate
$ws2.Range("H2").Value2 = $epoch
$ws2.Range("J2").Value2 = $isDependency
$ws2.Range("K2").Value2 = $dependencyFrom

# Row 3
$ws2.Range("A3").Value2 = $file2
$ws2.Range("B3").Value2 = $file1ext
$ws2.Range("C3").Value2 = $readyForHandoff
$ws2.Range("D3").Value2 = $target2
$ws2.Range("E3").Value2 = $zhDate
$ws2.Range("H3").Value2 = $epoch
$ws2.Range("J3").Value2 = $isDependency
$ws2.Range("K3").Value2 = $dependencyFrom

# Row 4
$ws2.Range("A4").Value2 = $file3
$ws2.Range("B4").Value2 = $file3ext
$ws2.Range("C4").Value2 = $readyForHandoff
$ws2.Range("D4").Value2 = $target3zh
$ws2.Range("E4").Value2 = $zhDate
$ws2.Range("H4").Value2 = $epoch
$ws2.Range("J4").Value2 = $include
$ws2.Range("K4").ClearContents()

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$file1", "", "", $file1)
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$target1", "", "", $target1)
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$file2", "", "", $file2)
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$target2", "", "", $target2)
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$file3", "", "", $file3)
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$target3zh", "", "", $target3zh)

# ===========================================================================
# Sheet 3: de-de
# ===========================================================================
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows(2).Copy()
$ws3.Rows(3).Insert(-4121)
$ws3.Rows(2).Copy()
$ws3.Rows(4).Insert(-4121)

# Row 2
$ws3.Range("A2").Value2 = $file1
$ws3.Range("B2").Value2 = $file1ext
$ws3.Range("C2").Value2 = $readyForHandoff
$ws3.Range("D2").Value2 = $target1
$ws3.Range("E2").Value2 = $overviewDate
$ws3.Range("H2").Value2 = $epoch
$ws3.Range("J2").Value2 = $isDependency
$ws3.Range("K2").Value2 = $dependencyFrom

# Row 3
$ws3.Range("A3").Value2 = $file2
$ws3.Range("B3").Value2 = $file1ext
$ws3.Range("C3").Value2 = $readyForHandoff
$ws3.Range("D3").Value2 = $target2
$ws3.Range("E3").Value2 = $overviewDate
$ws3.Range("H3").Value2 = $epoch
$ws3.Range("J3").Value2 = $isDependency
$ws3.Range("K3").Value2 = $dependencyFrom

# Row 4
$ws3.Range("A4").Value2 = $file3
$ws3.Range("B4").Value2 = $file3ext
$ws3.Range("C4").Value2 = $readyForHandoff
$ws3.Range("D4").Value2 = $target3de
$ws3.Range("E4").Value2 = $overviewDate
$ws3.Range("H4").Value2 = $epoch
$ws3.Range("J4").Value2 = $include
$ws3.Range("K4").ClearContents()

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$file1", "", "", $file1)
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$target1", "", "", $target1)
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$file2", "", "", $file2)
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$target2", "", "", $target2)
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/master/e2e/$file3", "", "", $file3)
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/master/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$target3de", "", "", $target3de)
